$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 464, shifting the existing rows 464:490 down to 467:493
$ws.Rows.Item(464).Insert()
$ws.Rows.Item(464).Insert()
$ws.Rows.Item(464).Insert()

# New row 464
$ws.Cells.Item(464,1).Value = 10
$ws.Cells.Item(464,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(464,3).Value = "La Araucanía"
$ws.Cells.Item(464,4).Value = 44753
$ws.Cells.Item(464,5).Value = 9
$ws.Cells.Item(464,6).Value = 100112032
$ws.Cells.Item(464,7).Value = "Zapallo italiano"
$ws.Cells.Item(464,8).Value = "Bola 8"
$ws.Cells.Item(464,9).Value = "Primera"
$ws.Cells.Item(464,10).Value = 80
$ws.Cells.Item(464,11).Value = 14000
$ws.Cells.Item(464,12).Value = 14000
$ws.Cells.Item(464,13).Value = 14000
$ws.Cells.Item(464,14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(464,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(464,16).Value = 233
$ws.Cells.Item(464,17).Value = 60
$ws.Cells.Item(464,18).Value = "Hortaliza"

# New row 465
$ws.Cells.Item(465,1).Value = 10
$ws.Cells.Item(465,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(465,3).Value = "La Araucanía"
$ws.Cells.Item(465,4).Value = 44753
$ws.Cells.Item(465,5).Value = 9
$ws.Cells.Item(465,6).Value = 100112032
$ws.Cells.Item(465,7).Value = "Zapallo italiano"
$ws.Cells.Item(465,8).Value = "Huracán"
$ws.Cells.Item(465,9).Value = "Primera"
$ws.Cells.Item(465,10).Value = 80
$ws.Cells.Item(465,11).Value = 14000
$ws.Cells.Item(465,12).Value = 14000
$ws.Cells.Item(465,13).Value = 14000
$ws.Cells.Item(465,14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(465,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(465,16).Value = 233
$ws.Cells.Item(465,17).Value = 60
$ws.Cells.Item(465,18).Value = "Hortaliza"

# New row 466
$ws.Cells.Item(466,1).Value = 10
$ws.Cells.Item(466,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(466,3).Value = "La Araucanía"
$ws.Cells.Item(466,4).Value = 44753
$ws.Cells.Item(466,5).Value = 9
$ws.Cells.Item(466,6).Value = 100112032
$ws.Cells.Item(466,7).Value = "Zapallo italiano"
$ws.Cells.Item(466,8).Value = "Sin especificar"
$ws.Cells.Item(466,9).Value = "Primera"
$ws.Cells.Item(466,10).Value = 350
$ws.Cells.Item(466,11).Value = 13000
$ws.Cells.Item(466,12).Value = 14000
$ws.Cells.Item(466,13).Value = 13429
$ws.Cells.Item(466,14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(466,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(466,16).Value = 224
$ws.Cells.Item(466,17).Value = 60
$ws.Cells.Item(466,18).Value = "Hortaliza"
